$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "DMWIRE-6"
$ws.Range("C6").Value = "40DC1FEF"

$ws.Range("C6").Select()
